$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.731.29'
$ws.Range("E2").Value = '  +1.48%  '

$ws.Range("D3").Value = '3.563.31'
$ws.Range("E3").Value = '  +1.37%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '''584.65'
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("D6").Value = '''188.34'
$ws.Range("E6").Value = '  +0.90%  '

$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").Value = '3.558.26'
$ws.Range("E7").Value = '  +1.42%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '''0.623'
$ws.Range("E8").Value = '  +1.26%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").Value = '''0.224'
$ws.Range("E10").Value = '  +10.07%  '

$ws.Range("D11").Value = '''0.648'
$ws.Range("E11").Value = '  -0.19%  '

$ws.Range("D12").Value = '''54.31'
$ws.Range("E12").Value = '  +0.60%  '

$ws.Range("D13").Value = '''0.0000316'
$ws.Range("E13").Value = '  +2.18%  '

$ws.Range("D14").Value = '''9.49'
$ws.Range("E14").Value = '  +0.01%  '

$ws.Range("D15").Value = '4.129.61'
$ws.Range("E15").Value = '  +1.35%  '

$ws.Range("D16").Value = '70.728.56'
$ws.Range("E16").Value = '  +1.53%  '

$ws.Range("D17").Value = '3.584.95'
$ws.Range("E17").Value = '  +2.15%  '

$ws.Range("D18").Value = '''19.13'
$ws.Range("E18").Value = '  -1.23%  '

$ws.Range("D19").Value = '''12.72'
$ws.Range("E19").Value = '  +3.15%  '

$ws.Range("D20").Value = '''568.50'
$ws.Range("E20").Value = '  +7.41%  '

$ws.Range("E21").Value = '  +0.73%  '

$ws.Range("E22").Value = '  -1.04%  '

$ws.Range("D23").Value = '''18.38'
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("D24").Value = '''4.59'
$ws.Range("E24").Value = '  +1.67%  '

$ws.Range("D25").Value = '''4.92'
$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("D26").Value = '''94.18'
$ws.Range("E26").Value = '  -0.36%  '

$ws.Range("D27").Value = '''11.14'
$ws.Range("E27").Value = '  -0.26%  '

$ws.Range("E28").Value = '  -1.23%  '

$ws.Range("D29").Value = '''9.29'
$ws.Range("E29").Value = '  +1.71%  '

$ws.Range("D30").Value = '''32.42'
$ws.Range("E30").Value = '  +1.64%  '

$ws.Range("D31").Value = '''7.14'
$ws.Range("E31").Value = '  -2.96%  '

$ws.Range("D32").Value = '''12.22'
$ws.Range("E32").Value = '  -3.29%  '

$ws.Range("E33").Value = '  +1.60%  '

$ws.Range("B34").Value = 'dogwifhat'
$ws.Range("C34").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D34").Value = '''3.84'
$ws.Range("E34").Value = '  +21.02%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '''63.38'
$ws.Range("E35").Value = '  -2.33%  '

$ws.Range("D36").Value = '''3.27'
$ws.Range("E36").Value = '  +6.09%  '

$ws.Range("B37").Value = 'TheGraph'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D37").Value = '''0.415'
$ws.Range("E37").Value = '  +0.53%  '

$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '''533.74'
$ws.Range("E38").Value = '  -1.20%  '

$ws.Range("D39").Value = '''38.51'
$ws.Range("E39").Value = '  +1.53%  '

$ws.Range("D40").Value = '''0.999'
$ws.Range("E40").Value = '  -0.10%  '

$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0797'
$ws.Range("E41").Value = '  +3.78%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '3.620.27'
$ws.Range("E42").Value = '  +8.46%  '

$ws.Range("D43").Value = '''0.139'
$ws.Range("E43").Value = '  +3.98%  '

$ws.Range("D44").Value = '''3.50'
$ws.Range("E44").Value = '  +3.56%  '

$ws.Range("D45").Value = '''0.0461'
$ws.Range("E45").Value = '  +4.18%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '''3.46'
$ws.Range("E46").Value = '  -0.02%  '

$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").Value = '''2.93'
$ws.Range("E47").Value = '  -2.01%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '''0.138'
$ws.Range("E48").Value = '  +2.40%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '''9.19'
$ws.Range("E49").Value = '  +2.83%  '

$ws.Range("D50").Value = '''0.998'
$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("E51").Value = '  +4.90%  '
